$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71. This shifts the existing rows 71..161
# down to 72..162, preserving all of their data (the diff shows every
# row from 71 onward effectively receiving the data that used to be in
# the row above it, with a brand-new record inserted at row 71 and the
# former last row (161) now duplicated/moved down to row 162).
$ws.Rows("71:71").Insert()

# Populate the newly inserted row 71 with the new record's data.
$ws.Range("A71").Value = 4
$ws.Range("B71").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C71").Value = "Los Lagos"
$ws.Range("D71").Value = 44413
$ws.Range("E71").Value = 10
$ws.Range("F71").Value = 100114014
$ws.Range("G71").Value = "Betarraga"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 1000
$ws.Range("L71").Value = 1000
$ws.Range("M71").Value = 1000
$ws.Range("N71").Value = "$/paquete 5 unidades"
$ws.Range("O71").Value = "Región del Maule"
$ws.Range("P71").Value = 200
$ws.Range("Q71").Value = 5
$ws.Range("R71").Value = "Hortaliza"
